$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: update title and link
$ws.Range("D12").Value = "iOS에서 빙 검색 위젯 설치하기"
$ws.Range("E12").Value = "https://tensorflow.blog/2023/06/19/ios%ec%97%90%ec%84%9c-%eb%b9%99-%ea%b2%80%ec%83%89-%ec%9c%84%ec%a0%af-%ec%84%a4%ec%b9%98%ed%95%98%ea%b8%b0/"

# Row 24: update title and link
$ws.Range("D24").Value = "Let's Verify Step by Step"
$ws.Range("E24").Value = "https://blog.naver.com/hist0134/223133196315"

# Row 51: update title and link
$ws.Range("D51").Value = "[vscode] Remote-SSH 확장 프로그램 활용하여 서버 ssh 접속하기"
$ws.Range("E51").Value = "https://bskyvision.com/entry/vscode-Remote-SSH-%ED%99%95%EC%9E%A5-%ED%94%84%EB%A1%9C%EA%B7%B8%EB%9E%A8-%ED%99%9C%EC%9A%A9%ED%95%98%EC%97%AC-%EC%84%9C%EB%B2%84-ssh-%EC%A0%91%EC%86%8D%ED%95%98%EA%B8%B0"
